$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Delete the now-unused "Comment/备注" column (H), which was always blank ---
$ws.Columns.Item(8).Delete()

# --- Re-label the header row from Chinese to the English panelized-export names ---
$ws.Range("A1").Value = "Designator"
$ws.Range("B1").Value = "Comment"
$ws.Range("C1").Value = "Footprint"
$ws.Range("D1").Value = "Mid X"
$ws.Range("E1").Value = "Mid Y"
$ws.Range("F1").Value = "Rotation"
$ws.Range("G1").Value = "Layer"
